$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 12.84410826922045
$ws.Range("D2").Value = 6.830014871696181
$ws.Range("E2").Value = 12.98835497247734
$ws.Range("F2").Value = 34.76132527280819
$ws.Range("G2").Value = 47.12024611927611
$ws.Range("H2").Value = 18.66121204132419
$ws.Range("K2").Value = 17.19094018507253
$ws.Range("L2").Value = 8.961393697488566
$ws.Range("M2").Value = 19.19045779502031
$ws.Range("N2").Value = 19.30386217300407
$ws.Range("C3").Value = 12.81451351674249
$ws.Range("D3").Value = 6.848036130100889
$ws.Range("E3").Value = 13.00623247203151
$ws.Range("F3").Value = 34.70824958172737
$ws.Range("G3").Value = 46.99469698285535
$ws.Range("H3").Value = 18.69788781148026
$ws.Range("K3").Value = 16.83706394540688
$ws.Range("L3").Value = 8.981874230174501
$ws.Range("M3").Value = 19.05759094293973
$ws.Range("N3").Value = 19.37716241867831
$ws.Range("C4").Value = 12.79926537814043
$ws.Range("D4").Value = 6.8595555708169
$ws.Range("E4").Value = 13.01946519715295
$ws.Range("F4").Value = 34.6860266774282
$ws.Range("G4").Value = 46.93345942157793
$ws.Range("H4").Value = 18.72449135834957
$ws.Range("K4").Value = 16.61994471517176
$ws.Range("L4").Value = 8.995295140787512
$ws.Range("M4").Value = 18.97960882837487
$ws.Range("N4").Value = 19.42413616522578
$ws.Range("C5").Value = 12.79379037634084
$ws.Range("D5").Value = 6.864364605313253
$ws.Range("E5").Value = 13.02542432961612
$ws.Range("F5").Value = 34.6795783242795
$ws.Range("G5").Value = 46.91249799524311
$ws.Range("H5").Value = 18.7363564929915
$ws.Range("K5").Value = 16.53163593031439
$ws.Range("L5").Value = 9.000977331796603
$ws.Range("M5").Value = 18.94876258076206
$ws.Range("N5").Value = 19.44377483809453
$ws.Range("C6").Value = 12.7929259681964
$ws.Range("D6").Value = 6.865170089078157
$ws.Range("E6").Value = 13.02644804697249
$ws.Range("F6").Value = 34.67866505688909
$ws.Range("G6").Value = 46.90925867836572
$ws.Range("H6").Value = 18.73838844006843
$ws.Range("K6").Value = 16.5169863316501
$ws.Range("L6").Value = 9.00193373494487
$ws.Range("M6").Value = 18.94369764465551
$ws.Range("N6").Value = 19.44706586574244
$ws.Range("C7").Value = 12.79918854488342
$ws.Range("D7").Value = 6.859619961766773
$ws.Range("E7").Value = 13.01954327042752
$ws.Range("F7").Value = 34.68592915500254
$ws.Range("G7").Value = 46.93316055332827
$ws.Range("H7").Value = 18.72464723401012
$ws.Range("K7").Value = 16.61875289620361
$ws.Range("L7").Value = 8.99537090962184
$ws.Range("M7").Value = 18.97918901635967
$ws.Range("N7").Value = 19.42439900655963
$ws.Range("C8").Value = 12.83329971842856
$ws.Range("D8").Value = 6.836134627163356
$ws.Range("E8").Value = 12.99405061779993
$ws.Range("F8").Value = 34.74087479708734
$ws.Range("G8").Value = 47.07367332163795
$ws.Range("H8").Value = 18.67300831118859
$ws.Range("K8").Value = 17.06896609941582
$ws.Range("L8").Value = 8.968280073090185
$ws.Range("M8").Value = 19.14391517429809
$ws.Range("N8").Value = 19.32872894326603
$ws.Range("C9").Value = 12.92319301523582
$ws.Range("D9").Value = 6.793660400635555
$ws.Range("E9").Value = 12.9619811236263
$ws.Range("F9").Value = 34.93070658381843
$ws.Range("G9").Value = 47.47439820788764
$ws.Range("H9").Value = 18.60428378243653
$ws.Range("K9").Value = 17.94764258307546
$ws.Range("L9").Value = 8.92184877085502
$ws.Range("M9").Value = 19.49423242963057
$ws.Range("N9").Value = 19.15664265262399
$ws.Range("C10").Value = 13.00295811912265
$ws.Range("D10").Value = 6.764602854553818
$ws.Range("E10").Value = 12.94937080053869
$ws.Range("F10").Value = 35.11982735347638
$ws.Range("G10").Value = 47.84399211043674
$ws.Range("H10").Value = 18.57379554575657
$ws.Range("K10").Value = 18.58361286181819
$ws.Range("L10").Value = 8.891792831520522
$ws.Range("M10").Value = 19.76646507452584
$ws.Range("N10").Value = 19.03955188001465
$ws.Range("C11").Value = 13.04214601217438
$ws.Range("D11").Value = 6.751842925465451
$ws.Range("E11").Value = 12.94601510858633
$ws.Range("F11").Value = 35.21650567133523
$ws.Range("G11").Value = 48.0281045459913
$ws.Range("H11").Value = 18.56429997017195
$ws.Range("K11").Value = 18.86946321568541
$ws.Range("L11").Value = 8.878995654585191
$ws.Range("M11").Value = 19.89313454884093
$ws.Range("N11").Value = 18.98828644419201
$ws.Range("C12").Value = 13.05739507660658
$ws.Range("D12").Value = 6.747076446359221
$ws.Range("E12").Value = 12.94508675363403
$ws.Range("F12").Value = 35.25463008704699
$ws.Range("G12").Value = 48.10008394037534
$ws.Range("H12").Value = 18.56133520492355
$ws.Range("K12").Value = 18.97710312840248
$ws.Range("L12").Value = 8.874275214121537
$ws.Range("M12").Value = 19.94147284895734
$ws.Range("N12").Value = 18.96915915941092
$ws.Range("C13").Value = 13.05409283930354
$ws.Range("D13").Value = 6.748100090892452
$ws.Range("E13").Value = 12.94527146556178
$ws.Range("F13").Value = 35.24635224729503
$ws.Range("G13").Value = 48.08448203452049
$ws.Range("H13").Value = 18.5619456261698
$ws.Range("K13").Value = 18.95394952956928
$ws.Range("L13").Value = 8.875286265458152
$ws.Range("M13").Value = 19.93104638842758
$ws.Range("N13").Value = 18.97326587792117
$ws.Range("C14").Value = 13.04339240315599
$ws.Range("D14").Value = 6.751449475928703
$ws.Range("E14").Value = 12.94593187099922
$ws.Range("F14").Value = 35.21961194214634
$ws.Range("G14").Value = 48.03398127599721
$ws.Range("H14").Value = 18.56404340006058
$ws.Range("K14").Value = 18.87833154716269
$ws.Range("L14").Value = 8.878604786209321
$ws.Range("M14").Value = 19.89710413367325
$ws.Range("N14").Value = 18.98670711256344
$ws.Range("C15").Value = 13.03689115572481
$ws.Range("D15").Value = 6.753509576931722
$ws.Range("E15").Value = 12.94638097358339
$ws.Range("F15").Value = 35.20342942331581
$ws.Range("G15").Value = 48.00334120732623
$ws.Range("H15").Value = 18.56541057898833
$ws.Range("K15").Value = 18.83193140079689
$ws.Range("L15").Value = 8.880653820105762
$ws.Range("M15").Value = 19.87636081468565
$ws.Range("N15").Value = 18.99497742807808
$ws.Range("C16").Value = 13.00045484053366
$ws.Range("D16").Value = 6.765445929038623
$ws.Range("E16").Value = 12.9496380152537
$ws.Range("F16").Value = 35.11372207912068
$ws.Range("G16").Value = 47.83227837696391
$ws.Range("H16").Value = 18.5745042973356
$ws.Range("K16").Value = 18.56485290645351
$ws.Range("L16").Value = 8.892646746056622
$ws.Range("M16").Value = 19.75824083652238
$ws.Range("N16").Value = 19.04294228803026
$ws.Range("C17").Value = 12.97884028891487
$ws.Range("D17").Value = 6.772885567791382
$ws.Range("E17").Value = 12.95224594065496
$ws.Range("F17").Value = 35.06140540675645
$ws.Range("G17").Value = 47.73140466416276
$ws.Range("H17").Value = 18.58120476875488
$ws.Range("K17").Value = 18.40004453195451
$ws.Range("L17").Value = 8.900228010324218
$ws.Range("M17").Value = 19.68647764493027
$ws.Range("N17").Value = 19.07287811917328
$ws.Range("C18").Value = 12.96668186719189
$ws.Range("D18").Value = 6.77720783671327
$ws.Range("E18").Value = 12.95397005087539
$ws.Range("F18").Value = 35.0323174798959
$ws.Range("G18").Value = 47.67489170653656
$ws.Range("H18").Value = 18.5854702678616
$ws.Range("K18").Value = 18.30493192375974
$ws.Range("L18").Value = 8.904670970055015
$ws.Range("M18").Value = 19.64546985007518
$ws.Range("N18").Value = 19.0902847462045
$ws.Range("C19").Value = 12.96261247221516
$ws.Range("D19").Value = 6.77867871517882
$ws.Range("E19").Value = 12.95459229173819
$ws.Range("F19").Value = 35.02264157453448
$ws.Range("G19").Value = 47.65601729294769
$ws.Range("H19").Value = 18.58698511133272
$ws.Range("K19").Value = 18.27267702131374
$ws.Range("L19").Value = 8.906189444664612
$ws.Range("M19").Value = 19.63163250332008
$ws.Range("N19").Value = 19.09621072972699
$ws.Range("C20").Value = 12.9811129234449
$ws.Range("D20").Value = 6.772089139700162
$ws.Range("E20").Value = 12.9519451284811
$ws.Range("F20").Value = 35.06687088585683
$ws.Range("G20").Value = 47.74198714319595
$ws.Range("H20").Value = 18.58044888092466
$ws.Range("K20").Value = 18.41762250670372
$ws.Range("L20").Value = 8.899412444245918
$ws.Range("M20").Value = 19.69408942681945
$ws.Range("N20").Value = 19.06967192024523
$ws.Range("C21").Value = 13.04652433405309
$ws.Range("D21").Value = 6.750463908609153
$ws.Range("E21").Value = 12.9457286025651
$ws.Range("F21").Value = 35.22742526095649
$ws.Range("G21").Value = 48.04875356513488
$ws.Range("H21").Value = 18.56341009267521
$ws.Range("K21").Value = 18.90055964385949
$ws.Range("L21").Value = 8.877626650671919
$ws.Range("M21").Value = 19.90706400520679
$ws.Range("N21").Value = 18.98275135277152
$ws.Range("C22").Value = 13.09165701762375
$ws.Range("D22").Value = 6.736711710187625
$ws.Range("E22").Value = 12.94366126486026
$ws.Range("F22").Value = 35.34117464699924
$ws.Range("G22").Value = 48.26239519301099
$ws.Range("H22").Value = 18.55595298522555
$ws.Range("K22").Value = 19.21261259555192
$ws.Range("L22").Value = 8.86412016831723
$ws.Range("M22").Value = 20.04840270359764
$ws.Range("N22").Value = 18.92760893327221
$ws.Range("C23").Value = 13.06735362121665
$ws.Range("D23").Value = 6.744016808008658
$ws.Range("E23").Value = 12.94458208162347
$ws.Range("F23").Value = 35.27966374385523
$ws.Range("G23").Value = 48.14718106948717
$ws.Range("H23").Value = 18.55959577018557
$ws.Range("K23").Value = 19.04642522809841
$ws.Range("L23").Value = 8.871261974575171
$ws.Range("M23").Value = 19.97278296081973
$ws.Range("N23").Value = 18.95688768667165
$ws.Range("C24").Value = 12.98008462979291
$ws.Range("D24").Value = 6.772449064236501
$ws.Range("E24").Value = 12.95208042551637
$ws.Range("F24").Value = 35.06439685893051
$ws.Range("G24").Value = 47.73719818939141
$ws.Range("H24").Value = 18.58078933052862
$ws.Range("K24").Value = 18.4096766237791
$ws.Range("L24").Value = 8.899780898780719
$ws.Range("M24").Value = 19.69064735936624
$ws.Range("N24").Value = 19.07112083173448
$ws.Range("C25").Value = 12.89644040744519
$ws.Range("D25").Value = 6.80477106765796
$ws.Range("E25").Value = 12.96873500899699
$ws.Range("F25").Value = 34.87059386196561
$ws.Range("G25").Value = 47.3526860550315
$ws.Range("H25").Value = 18.61937381436615
$ws.Range("K25").Value = 17.71113197356196
$ws.Range("L25").Value = 8.933695511560726
$ws.Range("M25").Value = 19.39672012446612
$ws.Range("N25").Value = 19.2015474032033
